$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 43301.52089262731

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newValue
}
